# Update automatico via Actualizar 02-09-2021 14-26-34
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new data rows 763-773 (2021-01-30 .. 2021-02-09) ---

# Copy the formatting of the last existing data row (762) down into the
# new rows so number formats / fonts / alignment match the rest of the
# "Dolar observado" column and the date column.
$ws.Range("A762:B762").Copy()
$ws.Range("A763:B773").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The rows that are missing a value in the source use the "--" placeholder
# style (bold, right aligned) instead of the numeric style - copy that
# formatting (taken from an existing "--" cell, B7) onto those cells.
$ws.Range("B7").Copy()
$ws.Range("B763:B764").PasteSpecial(-4122)
$ws.Range("B770:B771").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Dates (column A)
$ws.Range("A763").Value = 44226
$ws.Range("A764").Value = 44227
$ws.Range("A765").Value = 44228
$ws.Range("A766").Value = 44229
$ws.Range("A767").Value = 44230
$ws.Range("A768").Value = 44231
$ws.Range("A769").Value = 44232
$ws.Range("A770").Value = 44233
$ws.Range("A771").Value = 44234
$ws.Range("A772").Value = 44235
$ws.Range("A773").Value = 44236

# Values (column B) - "--" placeholders for missing/holiday dates,
# numeric observed dollar values otherwise.
$ws.Range("B763").Value = "--"
$ws.Range("B764").Value = "--"
$ws.Range("B765").Value = 734.62
$ws.Range("B766").Value = 731.66
$ws.Range("B767").Value = 734.86
$ws.Range("B768").Value = 730.53
$ws.Range("B769").Value = 737.23
$ws.Range("B770").Value = "--"
$ws.Range("B771").Value = "--"
$ws.Range("B772").Value = 736.65
$ws.Range("B773").Value = 735.07

# --- Column widths (narrower date column, wider value column) ---
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# --- Frozen pane / selection so the view matches the new bottom row ---
$ws.Range("B773").Select()

# --- Named ranges / filter database now cover the extended range ---
$wb.Names.Item("DOLAR_OBS_ADO!DOLAR_OBS_ADO").RefersTo = "=DOLAR_OBS_ADO!`$A`$1:`$B`$773"
$wb.Names.Item("DOLAR_OBS_ADO!_FilterDatabase").RefersTo = "=DOLAR_OBS_ADO!`$A`$3:`$B`$773"
